$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF)
# Copy the existing header formatting (bold font + border + centered
# alignment, same style index used by the other header cells) onto the
# two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Find the last used data row (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, "H").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value()  # column H = IP
    $ws.Cells.Item($r, 9).Value = 1           # column I = I0
    $ws.Cells.Item($r, 10).Value = $ipValue   # column J = IF
}
